# Updated capital structure database
# Applies updated values to rows 2 and 3 of the 'earnings_debt' sheet
# (both rows share identical data in this workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F3").Value = 0.15
$ws.Range("I2:I3").Value = 0.1488272921108742
$ws.Range("J2:J3").Value = 0.1259581301861937
$ws.Range("K2:K3").Value = 899
$ws.Range("L2:L3").Value = 0.09584221748400854
$ws.Range("M2:M3").Value = 709
$ws.Range("N2:N3").Value = 0.08583742947771134
$ws.Range("O2:O3").Value = 0.7886540600667408
$ws.Range("S2:S3").Value = 709
$ws.Range("U2:U3").Value = 7548
$ws.Range("V2:V3").Value = 0.9138235792634205
$ws.Range("W2:W3").Value = 0.06637135474344777
$ws.Range("X2:X3").Value = 0.07264538768372188
$ws.Range("Y2:Y3").Value = -0.006274032940274107
$ws.Range("Z2:Z3").Value = 0.8763897972531066
$ws.Range("AA2:AA3").Value = 0.1103884201762587
$ws.Range("AB2:AB3").Value = 0.06650216321179234
$ws.Range("AC2:AC3").Value = 0.04388625696446637
$ws.Range("AD2:AD3").Value = 1528
$ws.Range("AF2:AF3").Value = 1528
$ws.Range("AG2:AG3").Value = -6020
$ws.Range("AH2:AH3").Value = 0.1561127117431905
$ws.Range("AI2:AI3").Value = 0.08195226602306248
$ws.Range("AJ2:AJ3").Value = -2.687739976783642
$ws.Range("AK2:AK3").Value = -0.5424889609804452
$ws.Range("AL2:AL3").Value = 103
$ws.Range("AM2:AM3").Value = 103
$ws.Range("AO2:AO3").Value = 13.55339805825243
$ws.Range("AQ2:AQ3").Value = 13.55339805825243
